$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new bug row (row 14): bug description + status
$ws.Range("A14").Value = 'Apps with "" as name are just blank'
$ws.Range("B14").Value = "Yes they are"

# Update the selected cell to reflect the new active cell after the edit
$ws.Range("A15").Select()
